$d = $word.ActiveDocument

# Locate the paragraph that contains the "LOB1235: ..." heading text, then
# remove the trailing empty paragraph plus the two footer paragraphs
# ("Ver no Jupiter ..." and "(c) 2020 ...") that directly follow it, leaving
# the "LOB1235: ..." paragraph immediately followed by whatever came after
# that block.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "LOB1235:*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $delStart = $target.Range.End

    $p1 = $target.Next()
    $p2 = $p1.Next()
    $p3 = $p2.Next()
    $p4 = $p3.Next()
    $afterText = $p4

    $delEnd = $afterText.Range.Start

    $r = $d.Range($delStart, $delEnd)
    $r.Delete()
}
